$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the 2D "act" camera offset/rotation for the first scene row (villageScene)
$ws.Range("J2").Value = "0,4.2,5.5"
$ws.Range("K2").Value = "25,180"

# Bump the LoadingUI font size (column H, fontId used by H1/H4/H5)
$ws.Range("H1").Font.Size = 12
$ws.Range("H4").Font.Size = 12
$ws.Range("H5").Font.Size = 12

# Move the active selection
$ws.Range("K2").Select()
